# Apply cryptos list update (GitHub Actions refresh) to Sheet1.
# For each data row (2-51) set Coin, Link, Price and Volume(1h) to their
# updated values. Price cells are forced to Text format before the write
# so values like "66.00" / "1.901.58" are not reinterpreted as numbers,
# then the style is reset back to Normal so no extra formatting is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.862.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "

# Row 3: Ethereum
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "

# Row 4: TetherUSD
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5: BNB
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "

# Row 6: XRP
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.535"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.02%  "

# Row 7: USDC
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8: Cardano
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "

# Row 9: Dogecoin
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.53%  "

# Row 10: Solana
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.61%  "

# Row 11: TRON
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0895"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.66%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.901.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.36%  "

# Row 13: WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.678.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.84%  "

# Row 14: Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "

# Row 15: Polygon
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.44%  "

# Row 16: Litecoin
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.54%  "

# Row 17: WrappedBTC
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.874.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "

# Row 18: BitcoinCash
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.67%  "

# Row 19: Chainlink
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "

# Row 20: ShibaInu
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0731"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.12%  "

# Row 21: Dai
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22: Uniswap
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23: Toncoin
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.01%  "

# Row 24: Avalanche
$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.68%  "

# Row 25: Monero
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "

# Row 26: Cosmos
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.78%  "

# Row 27: Stellar
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.115"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.18%  "

# Row 28: EthereumClassic
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.06%  "

# Row 29: BinanceUSD
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30: Hedera
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0495"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.57%  "

# Row 31: PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "

# Row 32: Filecoin
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.09%  "

# Row 33: Maker
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.461.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.77%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.09%  "

# Row 35: LidoDAOToken
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.45%  "

# Row 36: HuobiToken
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "

# Row 37: ImmutableX
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "

# Row 38: ARBITRUM
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.898"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "

# Row 39: VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "

# Row 40: FraxShare
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.73%  "

# Row 41: PaxDollar
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42: MXToken
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "

# Row 43: WEMIXToken
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.979"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.58%  "

# Row 44: Aave
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "

# Row 45: RocketPoolETH
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.810.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.55%  "

# Row 46: TrustWalletToken
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.777"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "

# Row 47: Quant
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "

# Row 48: RenderToken
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

# Row 49: BabyDogeCoin
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "

# Row 50: Algorand
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.03%  "

# Row 51: Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0507"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.50%  "

